$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'35.543.64"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3
$ws.Range("D3").Value = "'1.911.04"
$ws.Range("E3").Value = "  +0.40%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "'0.703"
$ws.Range("E5").Value = "  +9.81%  "

# Row 6
$ws.Range("D6").Value = "'247.01"
$ws.Range("E6").Value = "  +0.32%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").Value = "'41.03"
$ws.Range("E8").Value = "  -2.08%  "

# Row 9
$ws.Range("E9").Value = "  +4.59%  "

# Row 10
$ws.Range("D10").Value = "'52.56"
$ws.Range("E10").Value = "  +7.73%  "

# Row 11
$ws.Range("D11").Value = "'0.0733"
$ws.Range("E11").Value = "  +3.99%  "

# Row 12
$ws.Range("D12").Value = "'0.0990"
$ws.Range("E12").Value = "  -0.86%  "

# Row 13
$ws.Range("D13").Value = "'2.186.42"
$ws.Range("E13").Value = "  +0.28%  "

# Row 14
$ws.Range("D14").Value = "'12.59"
$ws.Range("E14").Value = "  +1.65%  "

# Row 15
$ws.Range("D15").Value = "'0.717"
$ws.Range("E15").Value = "  +2.59%  "

# Row 16
$ws.Range("D16").Value = "'4.94"
$ws.Range("E16").Value = "  +2.74%  "

# Row 17
$ws.Range("D17").Value = "'1.911.93"
$ws.Range("E17").Value = "  +0.99%  "

# Row 18
$ws.Range("D18").Value = "'35.507.80"
$ws.Range("E18").Value = "  +0.07%  "

# Row 19
$ws.Range("D19").Value = "'73.39"
$ws.Range("E19").Value = "  +1.99%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0829"
$ws.Range("E20").Value = "  +0.02%  "

# Row 21
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'13.21"
$ws.Range("E21").Value = "  +4.37%  "

# Row 22
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'243.09"
$ws.Range("E22").Value = "  -0.01%  "

# Row 23
$ws.Range("D23").Value = "'5.07"
$ws.Range("E23").Value = "  +5.14%  "

# Row 24
$ws.Range("E24").Value = "  +0.04%  "

# Row 25
$ws.Range("E25").Value = "  +0.49%  "

# Row 26
$ws.Range("D26").Value = "'2.32"
$ws.Range("E26").Value = "  +3.65%  "

# Row 27
$ws.Range("D27").Value = "'169.98"
$ws.Range("E27").Value = "  -1.00%  "

# Row 28
$ws.Range("D28").Value = "'8.71"
$ws.Range("E28").Value = "  +1.94%  "

# Row 29
$ws.Range("D29").Value = "'18.87"
$ws.Range("E29").Value = "  +5.00%  "

# Row 30
$ws.Range("E30").Value = "  +5.15%  "

# Row 31
$ws.Range("D31").Value = "'4.109.46"
$ws.Range("E31").Value = "  +18.90%  "

# Row 32
$ws.Range("E32").Value = "  +3.26%  "

# Row 33
$ws.Range("D33").Value = "'0.0577"
$ws.Range("E33").Value = "  +1.32%  "

# Row 34
$ws.Range("E34").Value = "  +1.26%  "

# Row 35
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").Value = "'1.01"
$ws.Range("E35").Value = "  -0.04%  "

# Row 36
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'1.87"
$ws.Range("E36").Value = "  +5.67%  "

# Row 37
$ws.Range("D37").Value = "'0.916"
$ws.Range("E37").Value = "  -6.38%  "

# Row 38
$ws.Range("D38").Value = "'1.48"
$ws.Range("E38").Value = "  +10.75%  "

# Row 39
$ws.Range("E39").Value = "  -0.15%  "

# Row 40
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'17.24"
$ws.Range("E40").Value = "  +10.53%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'97.95"
$ws.Range("E41").Value = "  +6.85%  "

# Row 42
$ws.Range("E42").Value = "  +1.60%  "

# Row 43
$ws.Range("D43").Value = "'0.0210"
$ws.Range("E43").Value = "  +2.85%  "

# Row 44
$ws.Range("D44").Value = "'0.0651"
$ws.Range("E44").Value = "  +1.56%  "

# Row 45
$ws.Range("D45").Value = "'1.357.54"
$ws.Range("E45").Value = "  +0.78%  "

# Row 46
$ws.Range("D46").Value = "'2.43"
$ws.Range("E46").Value = "  +2.14%  "

# Row 47
$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D47").Value = "'2.42"
$ws.Range("E47").Value = "  +0.24%  "

# Row 48
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").Value = "'2.79"
$ws.Range("E48").Value = "  +1.20%  "

# Row 49
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "'46.00"
$ws.Range("E49").Value = "  -8.99%  "

# Row 50
$ws.Range("D50").Value = "'12.24"
$ws.Range("E50").Value = "  -5.03%  "

# Row 51
$ws.Range("D51").Value = "'6.57"
$ws.Range("E51").Value = "  -0.63%  "
